$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 19:48"

# Swap España / Colombia order: row 8 becomes España, row 9 becomes Colombia
$ws.Range("A8").Value = "España"
$ws.Range("A9").Value = "Colombia"

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-Row 4   8052784 14995 5199749 2632727 0 297 220308
Set-Row 5   7233670 60105 6293998 829084  0 694 110588
Set-Row 6   5105033 1625  4495269 458992  0 63  150772
Set-Row 8   925341  7118  0       0       0 80  33204
Set-Row 9   919083  0     798396  92702   0 0   27985
Set-Row 25  333770  2676  279100  44937   0 12  9733
Set-Row 27  296215  2184  243395  50780   0 19  2040
Set-Row 33  156946  3185  131462  22799   0 49  2685
Set-Row 65  55869   1245  24581   30809   0 13  479
Set-Row 72  44159   808   23364   18965   0 3   1830
Set-Row 110 9945    81    6347    3503    0 1   95
Set-Row 177 529     4     472     56      0 0   1
